$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.037.31"
$ws.Range("E2").Value = "  +1.01%  "

$ws.Range("D3").Value = "3.225.00"
$ws.Range("E3").Value = "  +0.32%  "

$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = "  -0.04%  "

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.52"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  -0.73%  "

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.57"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -7.18%  "

$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = "  -0.31%  "

$ws.Range("D8").Value = "3.224.64"
$ws.Range("E8").Value = "  +0.47%  "

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  -0.22%  "

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.163"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  -3.96%  "

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.29"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +6.90%  "

$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.479"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  +1.90%  "

$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000236"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  +0.98%  "

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.93"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  -5.03%  "

$ws.Range("D15").Value = "3.725.80"
$ws.Range("E15").Value = "  -0.06%  "

$ws.Range("D16").Value = "66.928.04"
$ws.Range("E16").Value = "  +1.18%  "

$ws.Range("D17").Value = "3.220.59"
$ws.Range("E17").Value = "  +0.36%  "

$ws.Range("E18").Value = "  -2.64%  "

$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.84"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  +0.55%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "502.98"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  -3.25%  "

$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.27"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -2.16%  "

$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.716"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  -3.65%  "

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.39"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -1.20%  "

$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.90"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  -1.55%  "

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.89"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  -1.40%  "

$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  +1.42%  "

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.08"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  -5.79%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "27.84"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  -1.91%  "

$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.04"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  -0.55%  "

$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.58"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  +2.95%  "

$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.17"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +5.05%  "

$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.52"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  +3.16%  "

$ws.Range("E33").Value = "  +0.10%  "

$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "498.77"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  -5.23%  "

$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.07"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  -2.44%  "

$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "54.41"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  +1.50%  "

$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.29"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  -5.02%  "

$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0414"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  -1.26%  "

$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0812"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  -1.58%  "

$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.52"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  -5.29%  "

$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.118"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +2.49%  "

$ws.Range("D42").Value = "2.843.98"
$ws.Range("E42").Value = "  +1.27%  "

$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.53"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  -4.04%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.251"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  +1.08%  "

$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  -0.10%  "

$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.01"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  -0.62%  "

$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.07"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  -1.54%  "

$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.02"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  +0.14%  "

$ws.Range("B49").Value = "PEPE"
$ws.Range("C49").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D49").Value = "0.0₃0529"
$ws.Range("E49").Value = "  -2.32%  "

$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.109"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  -1.26%  "

$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.10"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  -11.91%  "
